$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Stand"
$ws.Range("B1").Value = "Area"
$ws.Range("C1").Value = "Booked by"
$ws.Range("D1").Value = "Trade"
$ws.Range("E1").Value = "Time of booking"
$ws.Range("F1").Value = "Message to organizer"

# Row 2 (F2 left untouched - stays an empty cell, as in the diff)
$ws.Range("A2").Value = "1000-1002"
$ws.Range("B2").Value = "9x5m"
$ws.Range("C2").Value = "Leif Wallén"
$ws.Range("D2").Value = "Plats"
$ws.Range("E2").Value = "15-04-2013 13:48:05"

# Row 3
$ws.Range("A3").Value = "1713-1715"
$ws.Range("B3").Value = "9x5 m"
$ws.Range("C3").Value = "Grilltösen"
$ws.Range("D3").Value = "Gatuköksprodukter: olika sorters korv, hamburgare, pommes frites"
$ws.Range("E3").Value = "02-07-2013 08:16:28"
$ws.Range("F3").Value = "asdasdasd"

# Row 4
$ws.Range("A4").Value = "1227-12228"
$ws.Range("B4").Value = "6x5 m"
$ws.Range("C4").Value = "Nightmare on tour AB"
$ws.Range("D4").Value = "asdasd"
$ws.Range("E4").Value = "02-07-2013 08:16:37"
$ws.Range("F4").Value = "asdasd"

# Row 5
$ws.Range("A5").Value = "1113-1115"
$ws.Range("B5").Value = "9x5 m"
$ws.Range("C5").Value = "Marknadsmedia"
$ws.Range("D5").Value = "asdasd"
$ws.Range("E5").Value = "02-07-2013 08:16:47"
$ws.Range("F5").Value = "asdasd"
